$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.722364187240601
$ws.Range("B1").Value = 1.707981705665588
$ws.Range("C1").Value = 2.025704622268677
$ws.Range("D1").Value = 3.773809671401978
$ws.Range("E1").Value = 3.80626392364502
